$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (row 1) - copy the style (bold/centered) of the
# existing header cells onto the two new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M1").Value = "Future Dream"
$ws.Range("N1").Value = "Difficulty"

# Data rows 2-9 for "Future Dream" (M) and "Difficulty" (N) encoded lists
$futureDream = @(
    "[0, 1, 0, 0, 0, 0]",
    "[0, 0, 1, 0, 0, 0]",
    "[1, 0, 0, 0, 0, 0]",
    "[0, 0, 0, 0, 1, 0]",
    "[0, 0, 0, 1, 0, 0]",
    "[0, 1, 0, 0, 0, 0]",
    "[0, 0, 1, 0, 0, 0]",
    "[0, 0, 0, 0, 0, 1]"
)

$difficulty = @(
    "[0, 1, 0]",
    "[1, 0, 0]",
    "[0, 1, 0]",
    "[1, 0, 0]",
    "[1, 0, 0]",
    "[0, 1, 0]",
    "[1, 0, 0]",
    "[0, 0, 1]"
)

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $futureDream[$i]
    $ws.Cells.Item($row, 14).Value = $difficulty[$i]
}
